$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-2021")
$ws.Activate()

# ---------------------------------------------------------------------------
# Helper: style an "empty"/plain row cell the same way the existing plain
# (non-highlighted) rows in this tracker are styled - thin border + left
# aligned text, matching style index 24 used throughout the sheet.
# ---------------------------------------------------------------------------
function Format-Plain($rng) {
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4131
}

# New daily rows 15-22 (2021-01-14 .. 2021-01-21), continuing the tracker.
$rowsData = @(
    @{ Row = 15; No = 14; Date = 44210; Status = "Holiday" },
    @{ Row = 16; No = 15; Date = 44211; Status = "Holiday" },
    @{ Row = 17; No = 16; Date = 44212; Status = "Week off" },
    @{ Row = 18; No = 17; Date = 44213; Status = "Week off" }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $a = $ws.Cells.Item($row, 1)
    $a.Value = $r.No
    Format-Plain $a

    $b = $ws.Cells.Item($row, 2)
    $b.Value = $r.Date
    $b.NumberFormat = "[`$-14009]yyyy\-mm\-dd;@"
    Format-Plain $b

    $c = $ws.Cells.Item($row, 3)
    Format-Plain $c

    $d = $ws.Cells.Item($row, 4)
    $d.Value = $r.Status
    $d.Font.Bold = $true
    $d.Font.Color = 255
    $d.WrapText = $true
    $d.HorizontalAlignment = -4108
    $d.Borders.LineStyle = 1

    $e = $ws.Cells.Item($row, 5)
    Format-Plain $e

    $f = $ws.Cells.Item($row, 6)
    Format-Plain $f

    $g = $ws.Cells.Item($row, 7)
    Format-Plain $g
}

# Row 19 (2021-01-18): a completed task row with its own narrative comment.
$a19 = $ws.Cells.Item(19, 1)
$a19.Value = 18
Format-Plain $a19

$b19 = $ws.Cells.Item(19, 2)
$b19.Value = 44214
$b19.NumberFormat = "[`$-14009]yyyy\-mm\-dd;@"
Format-Plain $b19

$c19 = $ws.Cells.Item(19, 3)
$c19.Value = "Sonia"
$c19.WrapText = $true
$c19.HorizontalAlignment = -4131
$c19.Borders.LineStyle = 1

$commentText = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. `nRegression testing on Sonia Homac application`nRegression testing on Sonia Best Report ivc application`nRegression testing on Sonia Bic Report Soukastu application"
$d19 = $ws.Cells.Item(19, 4)
$d19.Value = $commentText
$d19.WrapText = $true
$d19.HorizontalAlignment = -4131
$d19.Borders.LineStyle = 1

$e19 = $ws.Cells.Item(19, 5)
$e19.Value = 1
$e19.NumberFormat = "0%"
Format-Plain $e19

$f19 = $ws.Cells.Item(19, 6)
$f19.Value = "Completed"
$f19.Interior.Color = 5287936
$f19.HorizontalAlignment = -4131
$f19.Borders.LineStyle = 1

$g19 = $ws.Cells.Item(19, 7)
Format-Plain $g19

$ws.Rows.Item(19).RowHeight = 60

# Rows 20-22 (2021-01-19 .. 2021-01-21): freshly added blank placeholder days.
$tailRows = @(
    @{ Row = 20; No = 19; Date = 44215 },
    @{ Row = 21; No = 20; Date = 44216 },
    @{ Row = 22; No = 21; Date = 44217 }
)

foreach ($r in $tailRows) {
    $row = $r.Row
    $a = $ws.Cells.Item($row, 1)
    $a.Value = $r.No
    Format-Plain $a

    $b = $ws.Cells.Item($row, 2)
    $b.Value = $r.Date
    $b.NumberFormat = "[`$-14009]yyyy\-mm\-dd;@"
    Format-Plain $b

    Format-Plain $ws.Cells.Item($row, 3)
    Format-Plain $ws.Cells.Item($row, 4)
    Format-Plain $ws.Cells.Item($row, 5)
    Format-Plain $ws.Cells.Item($row, 6)
    Format-Plain $ws.Cells.Item($row, 7)
}

# Update the view to match where the user ended up after entering the data.
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
